$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 120
$ws.Range("I10").Value = 120
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 120
$ws.Range("L10").Value = 120
$ws.Range("M10").Value = 173
$ws.Range("N10").Value = -706

$ws.Range("H20").Value = 9480.25
$ws.Range("I20").Value = 1307
$ws.Range("K20").Value = 1307
$ws.Range("M20").Value = -1077

$ws.Range("H21").Value = 21592.592
$ws.Range("I21").Value = 25419.834
$ws.Range("J21").Value = 16999.9
$ws.Range("K21").Value = 25419.834
$ws.Range("L21").Value = 16999.9
$ws.Range("M21").Value = -24951.834
$ws.Range("N21").Value = -17935.9

$ws.Range("H23").Value = 21592.592
$ws.Range("I23").Value = 25419.834
$ws.Range("J23").Value = 16999.9
$ws.Range("K23").Value = 25419.834
$ws.Range("L23").Value = 16999.9
$ws.Range("M23").Value = -25185.834
$ws.Range("N23").Value = -17467.9

$ws.Range("H35").Value = 9480.25
$ws.Range("I35").Value = 1307
$ws.Range("K35").Value = 1307
$ws.Range("M35").Value = -928

$ws.Range("H92").Value = 595.2727
$ws.Range("I92").Value = 574.8
$ws.Range("J92").Value = 800
$ws.Range("K92").Value = 574.8
$ws.Range("L92").Value = 800
$ws.Range("M92").Value = 673.2
$ws.Range("N92").Value = -3296

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8384.865
$ws.Range("I32").Value = 7208.487
$ws.Range("K32").Value = 7208.487
$ws.Range("M32").Value = -6921.487

$ws.Range("H37").Value = 27682.5
$ws.Range("I37").Value = 7500
$ws.Range("J37").Value = 47865
$ws.Range("K37").Value = 7500
$ws.Range("L37").Value = 47865
$ws.Range("M37").Value = -7227
$ws.Range("N37").Value = -48411

$ws.Range("H44").Value = 26995.428
$ws.Range("J44").Value = 26995.428
$ws.Range("L44").Value = 26995.428
$ws.Range("N44").Value = -27971.428

$ws.Range("H55").Value = 26498.75
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 26498.75
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 26498.75
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -27128.75

$ws.Range("H61").Value = 1845.5897
$ws.Range("I61").Value = 1518.2188
$ws.Range("J61").Value = 3342.1428
$ws.Range("K61").Value = 1518.2188
$ws.Range("L61").Value = 3342.1428
$ws.Range("M61").Value = -1306.2188
$ws.Range("N61").Value = -3766.1428

$ws.Range("H74").Value = 3372.16
$ws.Range("I74").Value = 3147.1428
$ws.Range("J74").Value = 4553.5
$ws.Range("K74").Value = 3147.1428
$ws.Range("L74").Value = 4553.5
$ws.Range("M74").Value = -2273.1428
$ws.Range("N74").Value = -6301.5

$ws.Range("H77").Value = 3372.16
$ws.Range("I77").Value = 3147.1428
$ws.Range("J77").Value = 4553.5
$ws.Range("K77").Value = 15735.714
$ws.Range("L77").Value = 22767.5
$ws.Range("M77").Value = -11367.714
$ws.Range("N77").Value = -31503.5

$ws.Range("H132").Value = 8476425
$ws.Range("I132").Value = 12501007
$ws.Range("J132").Value = 3621
$ws.Range("K132").Value = 37503021
$ws.Range("L132").Value = 10863
$ws.Range("M132").Value = -37500491
$ws.Range("N132").Value = -15923

$ws.Range("H136").Value = 1845.5897
$ws.Range("I136").Value = 1518.2188
$ws.Range("J136").Value = 3342.1428
$ws.Range("K136").Value = 4554.6564
$ws.Range("L136").Value = 10026.4284
$ws.Range("M136").Value = -2004.6564
$ws.Range("N136").Value = -15126.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1765.55
$ws.Range("I86").Value = 1813.6
$ws.Range("J86").Value = 1621.4
$ws.Range("K86").Value = 1813.6
$ws.Range("L86").Value = 1621.4
$ws.Range("M86").Value = -690.5999999999999
$ws.Range("N86").Value = -3867.4

$ws.Range("H89").Value = 1765.55
$ws.Range("I89").Value = 1813.6
$ws.Range("J89").Value = 1621.4
$ws.Range("K89").Value = 9068
$ws.Range("L89").Value = 8107
$ws.Range("M89").Value = -3452
$ws.Range("N89").Value = -19339

$ws.Range("H99").Value = 1844.3334
$ws.Range("I99").Value = 1747.6538
$ws.Range("K99").Value = 1747.6538
$ws.Range("M99").Value = -249.6538

$ws.Range("H105").Value = 2151.1052
$ws.Range("I105").Value = 2065.4546
$ws.Range("J105").Value = 2186
$ws.Range("K105").Value = 2065.4546
$ws.Range("L105").Value = 2186
$ws.Range("M105").Value = -318.4546
$ws.Range("N105").Value = -5680

$ws.Range("H134").Value = 2680.74
$ws.Range("I134").Value = 955.8889
$ws.Range("J134").Value = 3650.9688
$ws.Range("K134").Value = 2867.6667
$ws.Range("L134").Value = 10952.9064
$ws.Range("M134").Value = -332.6667000000002
$ws.Range("N134").Value = -16022.9064

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 39634.75
$ws.Range("J50").Value = 39634.75
$ws.Range("L50").Value = 39634.75
$ws.Range("N50").Value = -40884.75

$ws.Range("H60").Value = 24036.334
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 24036.334
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 24036.334
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -25058.334

$ws.Range("H62").Value = 2775.9167
$ws.Range("I62").Value = 2599
$ws.Range("J62").Value = 3023.6
$ws.Range("K62").Value = 2599
$ws.Range("L62").Value = 3023.6
$ws.Range("M62").Value = -1975
$ws.Range("N62").Value = -4271.6

$ws.Range("H65").Value = 2775.9167
$ws.Range("I65").Value = 2599
$ws.Range("J65").Value = 3023.6
$ws.Range("K65").Value = 12995
$ws.Range("L65").Value = 15118
$ws.Range("M65").Value = -9875
$ws.Range("N65").Value = -21358

$ws.Range("H74").Value = 13000
$ws.Range("J74").Value = 13000
$ws.Range("L74").Value = 13000
$ws.Range("N74").Value = -14748

$ws.Range("H77").Value = 13000
$ws.Range("J77").Value = 13000
$ws.Range("L77").Value = 39000
$ws.Range("N77").Value = -47736

$ws.Range("H99").Value = 2090.2
$ws.Range("I99").Value = 2163.25
$ws.Range("J99").Value = 2063.6365
$ws.Range("K99").Value = 2163.25
$ws.Range("L99").Value = 2063.6365
$ws.Range("M99").Value = -665.25
$ws.Range("N99").Value = -5059.636500000001

$ws.Range("H126").Value = 2090.2
$ws.Range("I126").Value = 2163.25
$ws.Range("J126").Value = 2063.6365
$ws.Range("K126").Value = 6489.75
$ws.Range("L126").Value = 6190.9095
$ws.Range("M126").Value = -4019.75
$ws.Range("N126").Value = -11130.9095

$ws.Range("H132").Value = 21871.564
$ws.Range("I132").Value = 1126.196
$ws.Range("J132").Value = 80650.11
$ws.Range("K132").Value = 3378.588
$ws.Range("L132").Value = 241950.33
$ws.Range("M132").Value = -848.5879999999997
$ws.Range("N132").Value = -247010.33

$ws.Range("H134").Value = 265635.1
$ws.Range("I134").Value = 1009.4722
$ws.Range("J134").Value = 826018.75
$ws.Range("K134").Value = 3028.4166
$ws.Range("L134").Value = 2478056.25
$ws.Range("M134").Value = -493.4166
$ws.Range("N134").Value = -2483126.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 5331.5
$ws.Range("J100").Value = 5331.5
$ws.Range("L100").Value = 15994.5
$ws.Range("N100").Value = -17616.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2246.195
$ws.Range("I7").Value = 1917.8214
$ws.Range("J7").Value = 2953.4614
$ws.Range("K7").Value = 1917.8214
$ws.Range("L7").Value = 2953.4614
$ws.Range("M7").Value = -1805.8214
$ws.Range("N7").Value = -3177.4614

$ws.Range("H122").Value = 45503.824
$ws.Range("I122").Value = 73234.86
$ws.Range("J122").Value = 2366.6667
$ws.Range("K122").Value = 219704.58
$ws.Range("L122").Value = 7100.000100000001
$ws.Range("M122").Value = -217254.58
$ws.Range("N122").Value = -12000.0001

$ws.Range("H126").Value = 2246.195
$ws.Range("I126").Value = 1917.8214
$ws.Range("J126").Value = 2953.4614
$ws.Range("K126").Value = 5753.4642
$ws.Range("L126").Value = 8860.3842
$ws.Range("M126").Value = -3283.4642
$ws.Range("N126").Value = -13800.3842

$ws.Range("H136").Value = 1894.2593
$ws.Range("I136").Value = 1711.1364
$ws.Range("K136").Value = 5133.4092
$ws.Range("M136").Value = -2583.4092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1581.8
$ws.Range("I132").Value = 1380.1212
$ws.Range("J132").Value = 2040.7931
$ws.Range("K132").Value = 4140.363600000001
$ws.Range("L132").Value = 6122.379300000001
$ws.Range("M132").Value = -1610.363600000001
$ws.Range("N132").Value = -11182.3793

$ws.Range("H136").Value = 11691.494
$ws.Range("I136").Value = 17140.85
$ws.Range("J136").Value = 1144.3549
$ws.Range("K136").Value = 51422.55
$ws.Range("L136").Value = 3433.0647
$ws.Range("M136").Value = -48872.55
$ws.Range("N136").Value = -8533.064699999999
Write-Output "done"
